# Update the ESPN comparison sheet with the next day's NBA matchups.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra row from the previous day's table (old table had 9 data rows, new has 8).
$ws.Range("A9:C9").Delete()

# Column A: date header + matchups, written top to bottom.
$ws.Range("A1").Value = "NBA, Thursday 7th Mar 2024"
$ws.Range("A2").Value = "Brooklyn Nets (25-37) vs Detroit Pistons (9-52)"
$ws.Range("A3").Value = "Minnesota Timberwolves (43-19) vs Indiana Pacers (35-29)"
$ws.Range("A4").Value = "Miami Heat (35-26) vs Dallas Mavericks (34-28)"
$ws.Range("A5").Value = "Toronto Raptors (23-39) vs Phoenix Suns (36-26)"
$ws.Range("A6").Value = "Boston Celtics (48-13) vs Denver Nuggets (42-20)"
$ws.Range("A7").Value = "Chicago Bulls (30-32) vs Golden State Warriors (33-28)"
$ws.Range("A8").Value = "San Antonio Spurs (13-49) vs Sacramento Kings (35-26)"

# Column B (Ballgorithm pick): rows 3-8 first, then row 2 (Detroit Pistons tie handled last).
$ws.Range("B3").Value = "Minnesota Timberwolves (74.19%)"
$ws.Range("B4").Value = "Miami Heat (56.67%)"
$ws.Range("B5").Value = "Phoenix Suns (60.61%)"
$ws.Range("B6").Value = "Boston Celtics (90.62%)"
$ws.Range("B7").Value = "Golden State Warriors (53.12%)"
$ws.Range("B8").Value = "Sacramento Kings (59.26%)"
$ws.Range("B2").Value = "Detroit Pistons (52.94%)"

# Column C (ESPN pick): row 2 right after, then rows 3-8.
$ws.Range("C2").Value = "Detroit Pistons (52.8%)"
$ws.Range("C3").Value = "Minnesota Timberwolves (55.9%)"
$ws.Range("C4").Value = "Dallas Mavericks (57.8%)"
$ws.Range("C5").Value = "Phoenix Suns (79.8%)"
$ws.Range("C6").Value = "Boston Celtics (69.4%)"
$ws.Range("C7").Value = "Golden State Warriors (72.8%)"
$ws.Range("C8").Value = "Sacramento Kings (76.5%)"

# Update the selection to match the saved view state.
$ws.Range("A8").Select()
